# Add outbreak info hits: three new sample rows (22-24) under columns A (Sample_ID)
# and J (Time length_(sampling)).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write new values in the same order the author typed them so new shared-string
# entries land in the expected order: L21_S1, L22N_S2, L22Q_S1, 17-23 June 2021,
# 24-30 June 2021.
$ws.Range("A22").Value = "L21_S1"
$ws.Range("A24").Value = "L22N_S2"
$ws.Range("A23").Value = "L22Q_S1"

$ws.Range("J22").Value = "17-23  June 2021"
$ws.Range("J23").Value = "24-30 June 2021"
$ws.Range("J24").Value = "24-30 June 2021"

# Match the saved selection/scroll state left behind in the workbook.
$ws.Range("L24").Select()
